# Commit: "deepika - 24th May"
# Edits applied to AutomationTestData.xlsx:
#  1. On the "D2HLeadSearch" sheet, fill the Password column (B) down for
#     every row in each Username "group" (previously only the first row of
#     each group had a Password, while the repeated blank rows below relied
#     on the now-removed "Visibility" column); then delete the entire
#     "Visibility" column (D), which collapses the sheet back to 3 columns.
#  2. Switch the active sheet from "D2HLeadSearch" to "D2HLeadCreation", and
#     update the remembered cell selections on both sheets.

$wb = $excel.ActiveWorkbook
$wsSearch = $wb.Worksheets.Item("D2HLeadSearch")

# --- Fill down column B (Password) for each grouped block of rows ---
# Group starting at row 2 (tushar.savdekar@sunrun.com / Sunrun+2000)
$wsSearch.Range("B3").Value = $wsSearch.Range("B2").Value()
$wsSearch.Range("B4").Value = $wsSearch.Range("B2").Value()
$wsSearch.Range("B5").Value = $wsSearch.Range("B2").Value()

# Group starting at row 6 (deepika.joshi@sunrun.com / Sunrun+1000)
$wsSearch.Range("B7").Value = $wsSearch.Range("B6").Value()
$wsSearch.Range("B8").Value = $wsSearch.Range("B6").Value()
$wsSearch.Range("B9").Value = $wsSearch.Range("B6").Value()

# Group starting at row 10 (manasi.kulkarni@sunrun.com / Sunrun+1000)
$wsSearch.Range("B11").Value = $wsSearch.Range("B10").Value()
$wsSearch.Range("B12").Value = $wsSearch.Range("B10").Value()
$wsSearch.Range("B13").Value = $wsSearch.Range("B10").Value()

# Group starting at row 14 (mandar35@yopmail.com / Sunrun+100)
$wsSearch.Range("B15").Value = $wsSearch.Range("B14").Value()
$wsSearch.Range("B16").Value = $wsSearch.Range("B14").Value()
$wsSearch.Range("B17").Value = $wsSearch.Range("B14").Value()

# --- Remove the "Visibility" column entirely ---
$wsSearch.Columns.Item(4).Delete()

# --- Update remembered selections / active sheet ---
[void]$wsSearch.Range("F11").Select()

$wsCreation = $wb.Worksheets.Item("D2HLeadCreation")
[void]$wsCreation.Activate()
[void]$wsCreation.Range("G16").Select()
